$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 192
$ws1.Range("F4").Value = 792
$ws1.Range("F6").Value = 10

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 192
$ws4.Range("F5").Value = 792
$ws4.Range("F7").Value = 10
